# Fix 0 hospitalisation error for India
# The "type" column (G) held a mix of data-type labels (date, timestamp,
# integer, nominal, ordinal, character) for what should uniformly be
# "character" entries. Unify every data row (G2:G42) to "character".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 7).Value = "character"
}

# Leave the unused type-label strings (timestamp/integer/nominal/ordinal)
# to fall out of the shared-string table naturally once nothing references
# them anymore.

# Scroll the view so row 21 is at the top (matches topLeftCell="C21") and
# select the corrected column so the selection highlights the fix.
$excel.Goto($ws.Range("C21"), $true) | Out-Null
$ws.Range("G2:G42").Select() | Out-Null
